$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1 ---
# Column B header relabeled from "field_number" to "site"
$ws.Range("B1").Value = "site"

# New headers for disease-related columns F:H, styled with font1 + bottom border (like A1:E1)
$ws.Range("F1").Value = "disease"
$ws.Range("G1").Value = "percent_diseased"
$ws.Range("H1").Value = "paling"
$ws.Range("F1:H1").Font.Color = 0
$ws.Range("F1:H1").Borders.Item(9).LineStyle = 1

# --- New data columns F (disease), G (percent_diseased), H (paling) for rows 2-34 ---
# Each inner array: row, disease, percent_diseased, paling
$diseaseData = @(
    @(2, "none", 0, 1),
    @(3, "none", 0, 1),
    @(4, "none", 0, 1),
    @(5, "none", 0, 1),
    @(6, "none", 0, 1),
    @(7, "none", 0, 0),
    @(8, "none", 0, 0),
    @(9, "none", 0, 0),
    @(10, "white_plague", 0.027272727, 1),
    @(11, "none", 0, 0),
    @(12, "none", 0, 1),
    @(13, "none", 0, 1),
    @(14, "none", 0, 1),
    @(15, "none", 0, 1),
    @(16, "none", 0, 0),
    @(17, "none", 0, 1),
    @(18, "none", 0, 1),
    @(19, "none", 0, 1),
    @(20, "none", 0, 0),
    @(21, "none", 0, 1),
    @(22, "none", 0, 0),
    @(23, "none", 0, 1),
    @(24, "none", 0, 0),
    @(25, "none", 0, 0),
    @(26, "none", 0, 0),
    @(27, "none", 0, 1),
    @(28, "none", 0, 0),
    @(29, "none", 0, 1),
    @(30, "none", 0, 0),
    @(31, "black_spot", 0.0125, 0),
    @(32, "black_spot", 0.003846154, 0),
    @(33, "none", 0, 1),
    @(34, "black_spot", 0.000666667, 1)
)

foreach ($row in $diseaseData) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
    $ws.Range("F" + $r + ":H" + $r).Font.Color = 0
}

# --- New blank rows 35-37, with F:H cells carrying the same data-row style (no values) ---
$ws.Range("F35:H37").Font.Color = 0

# --- Existing last row (38): keep A38/B38 text, add styled blank F38:H38 ---
$ws.Range("F38:H38").Font.Color = 0

# --- Update active cell selection ---
$ws.Range("J11").Select()
